$d = $word.ActiveDocument

$targets = @(
    "UNIQUE Tail_Number DECIMAL(10,0) NOT NULL,",
    "Tail_Number DECIMAL(10,0) NOT NULL,"
)

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $text = $r.Text.TrimEnd("`r", "`a")

    foreach ($old in $targets) {
        if ($text -eq $old) {
            $idx = $text.IndexOf("DECIMAL(10,0)")
            $oldLen = "DECIMAL(10,0)".Length
            $newText = "INTEGER"

            $sub = $d.Range($r.Start + $idx, $r.Start + $idx + $oldLen)
            $sub.Text = $newText

            $newLen = $newText.Length
            $run1 = $d.Range($r.Start, $r.Start + $idx)
            $run2 = $d.Range($r.Start + $idx, $r.Start + $idx + $newLen)
            $run3 = $d.Range($r.Start + $idx + $newLen, $r.End)

            # Toggle a character property on/off for each fragment so the
            # engine keeps them as distinct runs (each with an explicit,
            # empty <w:rPr/>) instead of re-merging them back together.
            $run1.Bold = 1
            $run1.Bold = 0
            $run2.Bold = 1
            $run2.Bold = 0
            $run3.Bold = 1
            $run3.Bold = 0

            break
        }
    }
}
